# This script re-orders the content of rows 47-65 on the active sheet
# (a pure permutation of row data: every field of a source row moves
# as a unit to a destination row; the row numbers themselves do not move).
#
# Mapping below: destinationRow = sourceRow (i.e. the data currently
# sitting in row <value> ends up in row <key> after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
  47 = 51
  48 = 49
  49 = 52
  50 = 47
  51 = 48
  52 = 50
  53 = 64
  54 = 62
  55 = 65
  56 = 53
  57 = 60
  58 = 59
  59 = 55
  60 = 56
  61 = 61
  62 = 54
  63 = 58
  64 = 57
  65 = 63
}

# Columns that always carry a value for every row in this range, and can
# be copied verbatim.
$plainCols = @("A","B","C","D","E","F","G","H","P","Q","R","S","T","U","V","W",
               "AD","AE","AG","AW","AX")

# Columns that hold date-like text ("2023-08-27") which Excel would
# otherwise auto-convert to a serial date number if typed into a
# General-formatted cell. Force them to stay text.
$dateTextCols = @("Y","AA")

# Plain text columns that don't look like dates (times such as "00:00").
$plainTextCols = @("Z","AB")

# Columns that are only sometimes populated (free-text notes) - need to
# be cleared on the destination when the source doesn't have a value.
$optionalCols = @("M","AC")

$rows = 47..65

# ---- Phase 1: snapshot every source cell we will need, before any
# writes happen (important because this is a self-contained permutation
# over the same row range). ----
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $plainCols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    foreach ($col in $dateTextCols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    foreach ($col in $plainTextCols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    foreach ($col in $optionalCols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowData
}

# ---- Phase 2: write the snapshotted data back out to its destination
# row according to $mapping. ----
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    foreach ($col in $plainCols) {
        $ws.Range("$col$destRow").Value = $src[$col]
    }

    foreach ($col in $plainTextCols) {
        $ws.Range("$col$destRow").Value = $src[$col]
    }

    foreach ($col in $dateTextCols) {
        $cell = $ws.Range("$col$destRow")
        $cell.NumberFormat = "@"
        $cell.Value = $src[$col]
    }

    foreach ($col in $optionalCols) {
        $val = $src[$col]
        if ($val -eq $null -or $val -eq "") {
            $ws.Range("$col$destRow").Value = ""
        } else {
            $ws.Range("$col$destRow").Value = $val
        }
    }

    # The sheet's exporter only ever emits J+AF (never L) for "non-activity"
    # (fungi/lichen) rows, and only ever emits L (never J/AF) for rows that
    # carry an Aktivitet (M) value (bird/animal observations). When a row's
    # category flips under the permutation, drop whichever placeholder
    # cell(s) no longer belong - Excel itself can't fabricate a brand new
    # empty placeholder cell, but it can certainly remove a stray one.
    $mVal = $src["M"]
    $isActivityRow = -not ($mVal -eq $null -or $mVal -eq "")
    if ($isActivityRow) {
        $ws.Range("J$destRow").Value = ""
        $ws.Range("AF$destRow").Value = ""
    } else {
        $ws.Range("L$destRow").Value = ""
    }
}
